$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 103
$ws.Range("H103").Value = 2765
$ws.Range("I103").Value = 600
$ws.Range("K103").Value = 1800
$ws.Range("M103").Value = -1214
# Row 105
$ws.Range("H105").Value = 34900
$ws.Range("J105").Value = 34900
$ws.Range("L105").Value = 34900
$ws.Range("N105").Value = -41888
# Row 138
$ws.Range("H138").Value = 4677.857
$ws.Range("I138").Value = 4165.6665
$ws.Range("K138").Value = 12496.9995
$ws.Range("M138").Value = -7356.999500000002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3791.3845
$ws.Range("I32").Value = 2856.75
$ws.Range("J32").Value = 15007
$ws.Range("K32").Value = 2856.75
$ws.Range("L32").Value = 15007
$ws.Range("M32").Value = -2569.75
$ws.Range("N32").Value = -15581
# Row 45
$ws.Range("H45").Value = 5093.4
$ws.Range("I45").Value = 1753.3334
$ws.Range("K45").Value = 1753.3334
$ws.Range("M45").Value = -1376.3334
# Row 74
$ws.Range("H74").Value = 12347248
$ws.Range("I74").Value = 13890341
$ws.Range("K74").Value = 13890341
$ws.Range("M74").Value = -13889467
# Row 77
$ws.Range("H77").Value = 12347248
$ws.Range("I77").Value = 13890341
$ws.Range("K77").Value = 69451705
$ws.Range("M77").Value = -69447337

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 29969.975
$ws.Range("I31").Value = 3687.3225
$ws.Range("K31").Value = 3687.3225
$ws.Range("M31").Value = -3392.3225
# Row 34
$ws.Range("H34").Value = 29969.975
$ws.Range("I34").Value = 3687.3225
$ws.Range("K34").Value = 3687.3225
$ws.Range("M34").Value = -3485.3225
# Row 55
$ws.Range("H55").Value = 19438.6
$ws.Range("I55").Value = 11098
$ws.Range("J55").Value = 24999
$ws.Range("K55").Value = 11098
$ws.Range("L55").Value = 24999
$ws.Range("M55").Value = -10783
$ws.Range("N55").Value = -25629
# Row 58
$ws.Range("H58").Value = 3052.879
$ws.Range("I58").Value = 1489.6666
$ws.Range("J58").Value = 10087.333
$ws.Range("K58").Value = 1489.6666
$ws.Range("L58").Value = 10087.333
$ws.Range("M58").Value = -1286.6666
$ws.Range("N58").Value = -10493.333
# Row 99
$ws.Range("H99").Value = 1955.5
$ws.Range("I99").Value = 1912
$ws.Range("K99").Value = 1912
$ws.Range("M99").Value = -414
# Row 122
$ws.Range("H122").Value = 9886.75
$ws.Range("I122").Value = 4298.875
$ws.Range("J122").Value = 21062.5
$ws.Range("K122").Value = 12896.625
$ws.Range("L122").Value = 63187.5
$ws.Range("M122").Value = -10446.625
$ws.Range("N122").Value = -68087.5
# Row 126
$ws.Range("H126").Value = 1955.5
$ws.Range("I126").Value = 1912
$ws.Range("K126").Value = 5736
$ws.Range("M126").Value = -3266
# Row 136
$ws.Range("H136").Value = 3052.879
$ws.Range("I136").Value = 1489.6666
$ws.Range("J136").Value = 10087.333
$ws.Range("K136").Value = 4468.9998
$ws.Range("L136").Value = 30261.999
$ws.Range("M136").Value = -1918.9998
$ws.Range("N136").Value = -35361.999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 6074.8
$ws.Range("J5").Value = 6500
$ws.Range("L5").Value = 19500
$ws.Range("N5").Value = -19724
# Row 40
$ws.Range("H40").Value = 62
$ws.Range("I40").Value = 13
$ws.Range("J40").Value = 76
$ws.Range("K40").Value = 52
$ws.Range("L40").Value = 304
$ws.Range("M40").Value = 17
$ws.Range("N40").Value = -442
# Row 88
$ws.Range("H88").Value = 17008
$ws.Range("I88").Value = 15000
$ws.Range("K88").Value = 45000
$ws.Range("M88").Value = -44572
# Row 91
$ws.Range("H91").Value = 17008
$ws.Range("I91").Value = 15000
$ws.Range("K91").Value = 45000
$ws.Range("M91").Value = -43518
# Row 92
$ws.Range("H92").Value = 1099.5
$ws.Range("J92").Value = 966
$ws.Range("L92").Value = 2898
$ws.Range("N92").Value = -5394
# Row 103
$ws.Range("H103").Value = 3664.2856
$ws.Range("I103").Value = 4000
$ws.Range("J103").Value = 3530
$ws.Range("K103").Value = 12000
$ws.Range("L103").Value = 10590
$ws.Range("M103").Value = -11121
$ws.Range("N103").Value = -12348
# Row 135
$ws.Range("H135").Value = 6074.8
$ws.Range("J135").Value = 6500
$ws.Range("L135").Value = 58500
$ws.Range("N135").Value = -63570
# Row 137
$ws.Range("H137").Value = 4496.5
$ws.Range("J137").Value = 6503.875
$ws.Range("L137").Value = 19511.625
$ws.Range("N137").Value = -29711.625
# Row 139
$ws.Range("H139").Value = 4756.643
$ws.Range("I139").Value = 1725.8889
$ws.Range("J139").Value = 10212
$ws.Range("K139").Value = 5177.6667
$ws.Range("L139").Value = 30636
$ws.Range("M139").Value = -37.66669999999976
$ws.Range("N139").Value = -40916
# Row 140
$ws.Range("H140").Value = 4364.45
$ws.Range("I140").Value = 2831.4375
$ws.Range("K140").Value = 8494.3125
$ws.Range("M140").Value = -3314.3125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 599
$ws.Range("I4").Value = 599
$ws.Range("K4").Value = 599
$ws.Range("M4").Value = -487
# Row 38
$ws.Range("H38").Value = 20000
$ws.Range("I38").Value = 20000
$ws.Range("K38").Value = 20000
$ws.Range("M38").Value = -19537
# Row 80
$ws.Range("H80").Value = 7257.6665
$ws.Range("J80").Value = 8388.6
$ws.Range("L80").Value = 8388.6
$ws.Range("N80").Value = -10384.6
# Row 83
$ws.Range("H83").Value = 7257.6665
$ws.Range("J83").Value = 8388.6
$ws.Range("L83").Value = 41943
$ws.Range("N83").Value = -51927
# Row 101
$ws.Range("H101").Value = 21976.924
$ws.Range("J101").Value = 21976.924
$ws.Range("L101").Value = 21976.924
$ws.Range("N101").Value = -28466.924
# Row 134
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6320.6665
$ws.Range("I22").Value = 1746.6
$ws.Range("J22").Value = 8079.923
$ws.Range("K22").Value = 1746.6
$ws.Range("L22").Value = 8079.923
$ws.Range("M22").Value = -1451.6
$ws.Range("N22").Value = -8669.922999999999
# Row 27
$ws.Range("H27").Value = 6320.6665
$ws.Range("I27").Value = 1746.6
$ws.Range("J27").Value = 8079.923
$ws.Range("K27").Value = 1746.6
$ws.Range("L27").Value = 8079.923
$ws.Range("M27").Value = -1639.6
$ws.Range("N27").Value = -8293.922999999999
# Row 68
$ws.Range("H68").Value = 4047.077
$ws.Range("I68").Value = 2408
$ws.Range("J68").Value = 5452
$ws.Range("K68").Value = 2408
$ws.Range("L68").Value = 5452
$ws.Range("M68").Value = -1659
$ws.Range("N68").Value = -6950
# Row 71
$ws.Range("H71").Value = 4047.077
$ws.Range("I71").Value = 2408
$ws.Range("J71").Value = 5452
$ws.Range("K71").Value = 12040
$ws.Range("L71").Value = 27260
$ws.Range("M71").Value = -8296
$ws.Range("N71").Value = -34748

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4633.1665
$ws.Range("I62").Value = 3600
$ws.Range("K62").Value = 3600
$ws.Range("M62").Value = -2976
# Row 65
$ws.Range("H65").Value = 4633.1665
$ws.Range("I65").Value = 3600
$ws.Range("K65").Value = 18000
$ws.Range("M65").Value = -14880
# Row 81
$ws.Range("H81").Value = 3826.5715
$ws.Range("I81").Value = 2659.2307
$ws.Range("K81").Value = 5318.4614
$ws.Range("M81").Value = -4257.4614
# Row 84
$ws.Range("H84").Value = 3826.5715
$ws.Range("I84").Value = 2659.2307
$ws.Range("K84").Value = 26592.307
$ws.Range("M84").Value = -21288.307

Write-Output "Applied 39 row updates across 8 sheets"
